{"js": "// Fix aliases with points in docx\n//\n// After the paragraph \"${col3}\" (inside the small one-column table that\n// also holds \"${col1}\" / \"${col2}\"), insert a new paragraph that renders\n// the dotted alias \"${col.nestedCol}\" as THREE separate runs -\n// \"${\", \"col.nestedCol\" and \"}\" - each carrying the same run properties\n// (w:lang=\"en-US\") as the surrounding text, mirroring how this template\n// already splits other dotted aliases (e.g. \"${Root.image}\") into\n// multiple runs instead of a single run.\n\n// Locate the anchor paragraph by its literal text so the script is not\n// dependent on hard-coded table/paragraph indices.\nconst anchorResults = context.document.body.search(\"${col3}\", { matchCase: true, matchWholeWord: false });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error(\"Could not find anchor paragraph containing '${col3}'\");\n}\n\nconst anchorParagraph = anchorResults.items[0].paragraphs.getFirst();\n\n// Create an empty paragraph right after it (still inside the same table\n// cell - Paragraph.insertParagraph keeps the new paragraph in the same\n// parent body/cell as the anchor).\nconst newParagraph = anchorParagraph.insertParagraph(\"\", \"After\");\n\n// Office.js' insertText()/insertParagraph(text, ...) merges the inserted\n// characters into a single run whenever the run formatting matches its\n// neighbour, which would collapse \"${\" + \"col.nestedCol\" + \"}\" into one\n// run. To reproduce the three discrete <w:r> runs from the diff, insert\n// raw OOXML (Flat OPC) for the paragraph instead - insertOoxml() lands\n// the runs verbatim, one <w:r> per literal in the fragment.\nconst flatOpcPackage = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>\\${</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>col.nestedCol</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>}</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nnewParagraph.getRange(\"Whole\").insertOoxml(flatOpcPackage, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Fix aliases with points in docx\n#\n# After the paragraph \"${col3}\" (inside the single-cell table that also\n# holds \"${col1}\" / \"${col2}\"), insert a new paragraph that renders the\n# dotted alias \"${col.nestedCol}\" as THREE separate runs - \"${\",\n# \"col.nestedCol\" and \"}\" - each carrying the same run properties\n# (w:lang=\"en-US\") as the surrounding text, mirroring how this template\n# already splits other dotted aliases (e.g. \"${Root.image}\") into\n# multiple runs instead of a single run.\n\n$d = $word.ActiveDocument\n\n# Locate the table/cell whose text contains the anchor \"${col3}\" instead\n# of relying on a hard-coded table index.\n$targetTable = $null\n$targetCell = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    for ($r = 1; $r -le $t.Rows.Count; $r++) {\n        for ($c = 1; $c -le $t.Columns.Count; $c++) {\n            $cell = $null\n            try {\n                $cell = $t.Cell($r, $c)\n            } catch {\n                $cell = $null\n            }\n            if ($cell -ne $null -and $cell.Range.Text -like \"*`${col3}*\") {\n                $targetTable = $t\n                $targetCell = $cell\n            }\n        }\n    }\n}\n\nif ($targetCell -eq $null) {\n    throw \"Could not find the table cell containing '`${col3}'\"\n}\n\n# The anchor paragraph is the cell's last paragraph (\"${col3}\").\n$paraCount = $targetCell.Range.Paragraphs.Count\n$anchorParagraph = $targetCell.Range.Paragraphs.Item($paraCount)\n\n# Insert a new, empty paragraph right after it (stays inside the cell).\n$anchorParagraph.Range.InsertParagraphAfter() | Out-Null\n\n# Re-fetch the cell/paragraph collection: InsertParagraphAfter() added a\n# block, so the new (now-last) paragraph is the one to fill in.\n$newParaCount = $targetCell.Range.Paragraphs.Count\n$newParagraph = $targetCell.Range.Paragraphs.Item($newParaCount)\n\n# Range.Text / Range.InsertAfter() coalesce same-formatted inserts into a\n# single run, which would collapse \"${\" + \"col.nestedCol\" + \"}\" into one\n# <w:r>. Use InsertXML (Flat OPC) instead so the three literals land as\n# three discrete runs, matching the target markup exactly.\n$flatOpcPackage = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>`${</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>col.nestedCol</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>}</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$newParagraph.Range.InsertXML($flatOpcPackage) | Out-Null\n"}
